# The document's single table has a trailing, entirely empty row
# (four empty cells matching the column widths 2269/3079/4717/5245)
# that was removed from the "Engagement Details" table. Remove the
# last row of the table to match.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$lastRow = $t.Rows.Item($t.Rows.Count)
$lastRow.Delete()
